$d = $word.ActiveDocument

# --- 1) Insert the new paragraph describing the effects of poor air quality ---
# Locate the heading paragraph "Efectos de la calidad del aire:" and insert a new,
# empty paragraph right after it (before the following empty NormalWeb paragraph),
# then fill that empty paragraph in with the full formatted content via its raw XML.
$found = $d.Content
[void]$found.Find.Execute("Efectos de la calidad del aire:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingIndex = $found.Paragraphs.Item(1).Index

$headingPara = $d.Paragraphs.Item($headingIndex)
$insertionRange = $headingPara.Range
$insertionRange.Collapse(0)
$insertionRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($headingIndex + 1)
$markRange = $newPara.Range

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-MX"/></w:rPr><w:t>La mala calidad del aire tiene varios efectos en los ciudadanos que se exponen a este fenómeno, como lo son ojos llorosos, tos, o ruido al respirar. Incluso si no se padece de enfermedades, la mala calidad del aire puede resultar en daños pulmonares y/o irritaciones al respirar. En caso de s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-MX"/></w:rPr><w:t>í</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-MX"/></w:rPr><w:t xml:space="preserve"> contar con una enfermedad previa, los efectos de la mala calidad del aire pueden resultar en un inmediato deterioro hacia la salud, especialmente si </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">dichos problemas están relacionados </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-MX"/></w:rPr><w:t xml:space="preserve">con </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr><w:t>dificultades</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-MX"/></w:rPr><w:t xml:space="preserve"> respiratorias o cardiovasculares. Suficiente exposición en situaciones específicas puede resultar mortal.</w:t></w:r></w:p>
'@

$markRange.InsertXML($newParagraphXml)

# --- 2) Move the lastRenderedPageBreak marker up one paragraph ---
# After the insertion above the document reflows, so the page break that used to
# render just before "El control de la cantidad del aire..." now falls one
# paragraph earlier, just before "Se ha buscado ir de estaciones...". Move the
# <w:lastRenderedPageBreak/> marker accordingly by rewriting both paragraphs' XML.

$withBreakXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3154BD48" w14:textId="77777777" w:rsidR="00C174F8" w:rsidRPr="00C174F8" w:rsidRDefault="00C174F8" w:rsidP="00C174F8"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00C174F8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr><w:lastRenderedPageBreak/><w:t>Se ha buscado ir de estaciones de mediciones fijas a sensores IoT, ya que las estaciones están limitadas a medir únicamente las concentraciones de ciertos puntos fijos en la ciudad y a la ves su coste de adquisición y mantenimiento elevados.</w:t></w:r></w:p>
'@

$withoutBreakXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6622EC6E" w14:textId="77777777" w:rsidR="00C174F8" w:rsidRPr="00C174F8" w:rsidRDefault="00C174F8" w:rsidP="00C174F8"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00C174F8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="es-ES"/></w:rPr><w:t>El control de la cantidad del aire se ha convertido un aspecto fundamental durante la pandemia, ya que es de gran relevancia conocer las condiciones del aire en varios sitios y encontrar patrones. Se observaron las siguientes tendencias:</w:t></w:r></w:p>
'@

$p1Found = $d.Content
[void]$p1Found.Find.Execute("Se ha buscado ir de estaciones", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p1Index = $p1Found.Paragraphs.Item(1).Index
$p1 = $d.Paragraphs.Item($p1Index)
$p1.Range.InsertXML($withBreakXml)

$p2Found = $d.Content
[void]$p2Found.Find.Execute("El control de la cantidad del aire", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2Index = $p2Found.Paragraphs.Item(1).Index
$p2 = $d.Paragraphs.Item($p2Index)
$p2.Range.InsertXML($withoutBreakXml)

Write-Output "Applied edits: new paragraph inserted, lastRenderedPageBreak moved"
